$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target team-name cell values after the shuffle (shared-string table was reshuffled,
# which changes the team-name text displayed in the data cells while the post/time
# header labels stay where they are).
$values = @{
    'A1' = '14:15-14:25'
    'A2' = 'Post 1'
    'B2' = 'fys5'
    'C2' = 'mat3'
    'D2' = 'hold2'
    'E2' = 'fys3'
    'A3' = 'Post 2'
    'B3' = 'dav'
    'C3' = 'mat2'
    'D3' = 'dat7'
    'E3' = 'hold1'
    'A4' = 'Post 3'
    'B4' = 'fys4'
    'C4' = 'fys2'
    'D4' = 'møk1'
    'E4' = 'it1'
    'A5' = 'Post 4'
    'B5' = 'nano'
    'C5' = 'it2'
    'D5' = 'TK'
    'E5' = 'mat1'
    'A6' = 'Post 5'
    'B6' = 'dat3'
    'C6' = 'møk2'
    'D6' = 'dat6'
    'E6' = 'dat4'
    'A7' = 'Pause'
    'B7' = 'fys1'
    'C7' = 'dat1'
    'D7' = 'dat5'
    'E7' = 'dat2'
    'A10' = '14:35-14:45'
    'A11' = 'Post 1'
    'B11' = 'fys2'
    'C11' = 'TK'
    'D11' = 'dat6'
    'E11' = 'dat2'
    'A12' = 'Post 2'
    'B12' = 'fys5'
    'C12' = 'møk1'
    'D12' = 'møk2'
    'E12' = 'dat5'
    'A13' = 'Post 3'
    'B13' = 'hold2'
    'C13' = 'dat7'
    'D13' = 'it2'
    'E13' = 'dat3'
    'A14' = 'Post 4'
    'B14' = 'mat3'
    'C14' = 'dav'
    'D14' = 'it1'
    'E14' = 'dat1'
    'A15' = 'Post 5'
    'B15' = 'fys3'
    'C15' = 'hold1'
    'D15' = 'fys4'
    'E15' = 'fys1'
    'A16' = 'Pause'
    'B16' = 'mat2'
    'C16' = 'nano'
    'D16' = 'mat1'
    'E16' = 'dat4'
    'A19' = '14:55-15:05'
    'A20' = 'Post 1'
    'B20' = 'hold1'
    'C20' = 'it1'
    'D20' = 'it2'
    'E20' = 'dat5'
    'A21' = 'Post 2'
    'B21' = 'mat3'
    'C21' = 'fys4'
    'D21' = 'mat1'
    'E21' = 'dat2'
    'A22' = 'Post 3'
    'B22' = 'fys5'
    'C22' = 'mat2'
    'D22' = 'nano'
    'E22' = 'dat4'
    'A23' = 'Post 4'
    'B23' = 'hold2'
    'C23' = 'fys2'
    'D23' = 'møk2'
    'E23' = 'fys1'
    'A24' = 'Post 5'
    'B24' = 'dat7'
    'C24' = 'møk1'
    'D24' = 'TK'
    'E24' = 'dat1'
    'A25' = 'Pause'
    'B25' = 'fys3'
    'C25' = 'dav'
    'D25' = 'dat3'
    'E25' = 'dat6'
    'A28' = '15:15-15:25'
    'A29' = 'Post 1'
    'B29' = 'dav'
    'C29' = 'møk1'
    'D29' = 'nano'
    'E29' = 'dat3'
    'A30' = 'Post 2'
    'B30' = 'fys3'
    'C30' = 'fys2'
    'D30' = 'it2'
    'E30' = 'dat4'
    'A31' = 'Post 3'
    'B31' = 'fys1'
    'C31' = 'dat1'
    'D31' = 'dat5'
    'E31' = 'dat2'
    'A32' = 'Post 4'
    'B32' = 'fys5'
    'C32' = 'dat7'
    'D32' = 'fys4'
    'E32' = 'dat6'
    'A33' = 'Post 5'
    'B33' = 'hold2'
    'C33' = 'mat2'
    'D33' = 'it1'
    'E33' = 'mat1'
    'A34' = 'Pause'
    'B34' = 'mat3'
    'C34' = 'hold1'
    'D34' = 'TK'
    'E34' = 'møk2'
    'A37' = '15:35-15:45'
    'A38' = 'Post 1'
    'B38' = 'dat7'
    'C38' = 'mat1'
    'D38' = 'dat4'
    'E38' = 'fys1'
    'A39' = 'Post 2'
    'B39' = 'hold2'
    'C39' = 'nano'
    'D39' = 'dat6'
    'E39' = 'dat1'
    'A40' = 'Post 3'
    'B40' = 'mat3'
    'C40' = 'hold1'
    'D40' = 'TK'
    'E40' = 'møk2'
    'A41' = 'Post 4'
    'B41' = 'fys3'
    'C41' = 'mat2'
    'D41' = 'dat3'
    'E41' = 'dat5'
    'A42' = 'Post 5'
    'B42' = 'fys5'
    'C42' = 'dav'
    'D42' = 'it2'
    'E42' = 'dat2'
    'A43' = 'Pause'
    'B43' = 'fys4'
    'C43' = 'fys2'
    'D43' = 'møk1'
    'E43' = 'it1'
    'A46' = '15:55-16:05'
    'A47' = 'Post 1'
    'B47' = 'mat2'
    'C47' = 'fys4'
    'D47' = 'møk2'
    'E47' = 'dat1'
    'A48' = 'Post 2'
    'B48' = 'it1'
    'C48' = 'TK'
    'D48' = 'dat3'
    'E48' = 'fys1'
    'A49' = 'Post 3'
    'B49' = 'fys3'
    'C49' = 'dav'
    'D49' = 'mat1'
    'E49' = 'dat6'
    'A50' = 'Post 4'
    'B50' = 'hold1'
    'C50' = 'møk1'
    'D50' = 'dat4'
    'E50' = 'dat2'
    'A51' = 'Post 5'
    'B51' = 'mat3'
    'C51' = 'fys2'
    'D51' = 'nano'
    'E51' = 'dat5'
    'A52' = 'Pause'
    'B52' = 'fys5'
    'C52' = 'hold2'
    'D52' = 'dat7'
    'E52' = 'it2'
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
